$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: dimension/measure identifier relabeling
$ws.Range("A3").Value = "iaest-dimension:edad-grandes-grupos"
$ws.Range("E3").Value = "sdmx-dimension:refArea"
$ws.Range("F3").Value = "sdmx-dimension:refArea"
$ws.Range("G3").Value = "sdmx-dimension:refArea"
$ws.Range("I3").Value = "iaest-dimension:sexo"
$ws.Range("J3").Value = "iaest-dimension:relacion-lugar-de-residencia-y-nacimiento"

# Row 4: "medida" -> "dim" for the columns that became dimensions
$ws.Range("A4").Value = "dim"
$ws.Range("E4").Value = "dim"
$ws.Range("F4").Value = "dim"
$ws.Range("G4").Value = "dim"
$ws.Range("I4").Value = "dim"
$ws.Range("J4").Value = "dim"

# Row 5: data type updates
$ws.Range("A5").Value = "skos:Concept"
$ws.Range("E5").Value = "URI-Comunidad"
$ws.Range("F5").Value = "URI-comarca"
$ws.Range("G5").Value = "URI-Provincia"
$ws.Range("I5").Value = "skos:Concept"
$ws.Range("J5").Value = "skos:Concept"

# Row 6: new mapping file references (copy formatting from row 1 so the
# new cells pick up the same style used throughout the sheet)
$ws.Range("A1").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("J6").PasteSpecial(-4122)

$ws.Range("A6").Value = "mapping-edad-grandes-grupos.xlsx"
$ws.Range("I6").Value = "mapping-sexo.xlsx"
$ws.Range("J6").Value = "mapping-relacion-lugar-de-residencia-y-nacimiento.xlsx"

$excel.CutCopyMode = $false
